# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker data table (rows 16-29, columns B:G) is re-grouped so that all
# periods for "JOSE MADERA PEREZ" (73095916) come first, followed by all
# periods for "MARTHA LUCIA RENGIFO SANGUINO" (45437880), each ordered by
# descending period (2108 -> 2102). The Valor Mora / Salario Basico figures
# travel together with their (worker, period) pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2108", 33942, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2107", 37570, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2106", 37570, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2105", 37570, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2104", 37570, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2103", 37570, 939249),
    @("CC", "73095916", "JOSE MADERA PEREZ",             "2102", 37570, 939249),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2108", 33942, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2107", 35112, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2106", 35112, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2105", 35112, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2104", 36341, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2103", 36341, 908526),
    @("CC", "45437880", "MARTHA LUCIA RENGIFO SANGUINO", "2102", 36341, 908526)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G: Salario Basico
}
